# Daily attendance processing - reorder "Recorded By" (column G) entries.
# For every row where the comma-separated list of recorders ends with the
# token "System" (exact case), reverse the order of the comma-separated
# tokens so "System" appears first. Lists that do not end in "System" are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7 ("Recorded By")
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ','
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }

        if ($parts.Length -gt 0 -and $parts[$parts.Length - 1] -eq 'System') {
            $n = $parts.Length
            $reversed = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $cell.Value = [string]::Join(', ', $reversed)
        }
    }
}
